$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new timesheet entry as row 9, matching the formatting of the
# existing data row above it (date style, hours style, wrapped description).
$ws.Rows.Item(8).Copy()
$ws.Rows.Item(9).PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = $false

$ws.Range("A9").Value = 41987
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = "Styling for header, main menu, mobile responsive styling for global nav and main nav, styles for page text like page headings, standfirst, secondary headings and paragraphs."

# Taller row to accommodate the wrapped description text (matches the
# height used by other multi-line description rows in the sheet).
$ws.Rows.Item(9).RowHeight = 30

# Move/refresh the active selection like the authored workbook.
$ws.Range("C11").Select()
